$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: helper label (source file name) ---
$ws.Range("B2").Value = "HourProcessor.cs"

# --- Row 5 / Row 6 helper labels (re-asserted, unchanged text) ---
$ws.Range("N5").Value = "actual duty needed"
$ws.Range("Q6").Value = "overtime"

# --- Row 7: "1" header factor becomes an annotated "1 = 1" label,
#     and the existing "x * y" factor labels gain their computed "= result" ---
$ws.Range("C7").Value = "1 = 1"
$ws.Range("G7").Value = "1.3 * 1.1 = 1.43"
$ws.Range("H7").Value = "2 * 1.1 = 2.2"

# --- Row 8: nsu_overtime_* labels (re-asserted, unchanged text) ---
$ws.Range("C8").Value = "nsu_overtime_day_normal"
$ws.Range("D8").Value = "nsu_overtime_day_special"
$ws.Range("E8").Value = "nsu_overtime_day_regular"
$ws.Range("F8").Value = "nsu_overtime_night_normal"
$ws.Range("G8").Value = "nsu_overtime_night_special"
$ws.Range("H8").Value = "nsu_overtime_night_regular"

# --- Row 9: pay-rate formula labels now show the computed result ---
$ws.Range("D9").Value = "1.3 * 1.3 = 1.69"
$ws.Range("E9").Value = "2 * 1.3 = 2.6"
$ws.Range("F9").Value = "1.1* 1.25 = 1.375"
$ws.Range("G9").Value = "1.3 * 1.1 * 1.3 = 1.859"
$ws.Range("H9").Value = "2 * 1.1 * 1.3 = 2.86"

# --- Row 11: pay-rate formula labels now show the computed result ---
$ws.Range("F11").Value = "1.3 * 1.1 = 1.43"
$ws.Range("G11").Value = "1.5 * 1.1 = 1.65"
$ws.Range("H11").Value = "2.6 * 1.1 = 2.86"

# --- Row 12: sun_overtime_* labels (re-asserted, unchanged text) ---
$ws.Range("C12").Value = "sun_overtime_day_normal"
$ws.Range("D12").Value = "sun_overtime_day_special"
$ws.Range("E12").Value = "sun_overtime_day_regular"
$ws.Range("F12").Value = "sun_overtime_night_normal"
$ws.Range("G12").Value = "sun_overtime_night_special"
$ws.Range("H12").Value = "sun_overtime_night_regular"

# --- Row 13: pay-rate formula labels now show the computed result ---
$ws.Range("C13").Value = "1.3 * 1.25 = 1.625"
$ws.Range("D13").Value = "1.5 * 1.3= 1.95"
$ws.Range("E13").Value = "2.6 * 1.3 = 3.38"
$ws.Range("F13").Value = "1.3 * 1.1 * 1.3 = 1.859"
$ws.Range("G13").Value = "1.5 * 1.1 * 1.3 = 2.145"
$ws.Range("H13").Value = "2.6 * 1.1 * 1.3 = 3.718"

# --- Refresh the view: clear the stale scroll/selection state left over
#     from editing, matching the author's final cursor position ---
$ws.Range("F23").Select()
